$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.722.44'
$ws.Range("E2").Value = '  -1.10%  '
$ws.Range("D3").Value = '1.601.14'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.514'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  -1.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").Value = '1.820.19'
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("D13").Value = '1.596.40'
$ws.Range("E13").Value = '  -1.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.524'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("D17").Value = '26.686.31'
$ws.Range("E17").Value = '  -1.22%  '
$ws.Range("E18").Value = '  -0.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '210.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.31%  '
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.66'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.48%  '
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.116'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.36'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.91%  '
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("E32").Value = '  -2.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.672'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.90'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = '1.299.23'
$ws.Range("E35").Value = '  -2.68%  '
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("E37").Value = '  -5.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0172'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.844'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.69%  '
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.790'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.36'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.98'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("D45").Value = '1.734.12'
$ws.Range("E45").Value = '  -1.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.899'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("E49").Value = '  -1.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0990'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0504'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.06%  '
